# Cleaned up TrialDetailViewLoad test
#
# - Adds a new "TrialDetailView" worksheet (after DynamicListingPage) that
#   captures three sample trial-detail-view hits (Advanced / Basic / Custom
#   content types), mirroring the existing Path/ContentType-style load
#   sheets.
# - Leaves the previously-active "DynamicListingPage" sheet selected on the
#   A1:B1 header range instead of the stray A8 cell, and makes the new
#   sheet the active tab.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$dynamicListing = $wb.Worksheets.Item("DynamicListingPage")

# --- New sheet, appended after the last existing tab -----------------
$trial = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$trial.Name = "TrialDetailView"

# --- Header row (copy formatting from another sheet's bold header so we
#     reuse the existing shared header style rather than minting a new one)
$trial.Range("A1").Value = "Path"
$trial.Range("B1").Value = "ContentType"
$ws1.Range("A1").Copy()
$trial.Range("A1:B1").PasteSpecial(-4122)

# --- Sample data rows --------------------------------------------------
$trial.Range("A2").Value = "?t=C4911&q=nivolumab&loc=0&tid=S1609&rl=2&id=NCI-2016-01041&pn=1&ni=10"
$trial.Range("B2").Value = "Advanced"

$trial.Range("A3").Value = "?q=ipilimumab&loc=1&z=20850&zp=100&rl=1&id=NCI-2016-01041&pn=1&ni=10"
$trial.Range("B3").Value = "Basic"

$trial.Range("A4").Value = "?id=NCI-2016-01041&r=1"
$trial.Range("B4").Value = "Custom"

# --- Column widths (best-fit-ish, matching the look of the other sheets)
$trial.Columns.Item(1).ColumnWidth = 73.333333
$trial.Columns.Item(2).ColumnWidth = 15.333333

# --- Selection / active-tab bookkeeping --------------------------------
$dynamicListing.Range("A1:B1").Select()

$trial.Activate()
$trial.Range("A5").Select()

Write-Output "TrialDetailView sheet added"
